# Weekly data refresh: insert two new price records (rows 187-188) for
# "Feria Lagunitas de Puerto Montt" / Platano, pushing the existing rows
# down by two (old row 187 -> new row 189, ... old row 286 -> new row 288).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 187 (shifts everything below down by one,
# twice in a row).
$ws.Rows.Item(187).Insert()
$ws.Rows.Item(187).Insert()

# New row 187: Barraganete / Primera
$ws.Cells.Item(187, 1).Value = 4
$ws.Cells.Item(187, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(187, 3).Value = "Los Lagos"
$ws.Cells.Item(187, 4).Value = 44523
$ws.Cells.Item(187, 5).Value = 10
$ws.Cells.Item(187, 6).Value = "Fruta"
$ws.Cells.Item(187, 7).Value = 100108
$ws.Cells.Item(187, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(187, 9).Value = 100108006
$ws.Cells.Item(187, 10).Value = "Plátano"
$ws.Cells.Item(187, 11).Value = "Barraganete"
$ws.Cells.Item(187, 12).Value = "Primera"
$ws.Cells.Item(187, 13).Value = 200
$ws.Cells.Item(187, 14).Value = 32000
$ws.Cells.Item(187, 15).Value = 33000
$ws.Cells.Item(187, 16).Value = 32500
$ws.Cells.Item(187, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(187, 18).Value = "Ecuador"
$ws.Cells.Item(187, 19).Value = 1625
$ws.Cells.Item(187, 20).Value = 20

# New row 188: Sin especificar / Primera Pintón
$ws.Cells.Item(188, 1).Value = 4
$ws.Cells.Item(188, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(188, 3).Value = "Los Lagos"
$ws.Cells.Item(188, 4).Value = 44523
$ws.Cells.Item(188, 5).Value = 10
$ws.Cells.Item(188, 6).Value = "Fruta"
$ws.Cells.Item(188, 7).Value = 100108
$ws.Cells.Item(188, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(188, 9).Value = 100108006
$ws.Cells.Item(188, 10).Value = "Plátano"
$ws.Cells.Item(188, 11).Value = "Sin especificar"
$ws.Cells.Item(188, 12).Value = "Primera Pintón"
$ws.Cells.Item(188, 13).Value = 1400
$ws.Cells.Item(188, 14).Value = 23000
$ws.Cells.Item(188, 15).Value = 24000
$ws.Cells.Item(188, 16).Value = 23500
$ws.Cells.Item(188, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(188, 18).Value = "Ecuador"
$ws.Cells.Item(188, 19).Value = 1175
$ws.Cells.Item(188, 20).Value = 20
